# Updated cryptos list values (prices stored as text to match original inlineStr cells)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.506.66"
$ws.Range('E2').Value = '  +1.58%  '

$ws.Range('D3').Value = "'1.676.38"
$ws.Range('E3').Value = '  +1.75%  '

$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = "'219.71"
$ws.Range('E5').Value = '  +1.34%  '

$ws.Range('D6').Value = "'0.5315"

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('E8').Value = '  +3.19%  '

$ws.Range('D9').Value = "'0.06396"
$ws.Range('E9').Value = '  -0.67%  '

$ws.Range('D10').Value = "'21.77"
$ws.Range('E10').Value = '  +4.35%  '

$ws.Range('D11').Value = "'0.07794"
$ws.Range('E11').Value = '  +1.19%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.504"
$ws.Range('E12').Value = '  +1.75%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.669.66"
$ws.Range('E13').Value = '  +1.24%  '

$ws.Range('D14').Value = "'0.5584"

$ws.Range('D15').Value = "'0.0₅8353"
$ws.Range('E15').Value = '  +0.74%  '

$ws.Range('D16').Value = "'65.69"
$ws.Range('E16').Value = '  +0.74%  '

$ws.Range('D17').Value = "'26.522.46"
$ws.Range('E17').Value = '  +1.62%  '

$ws.Range('E18').Value = '  -0.09%  '

$ws.Range('D19').Value = "'4.785"
$ws.Range('E19').Value = '  +1.01%  '

$ws.Range('D20').Value = "'192.66"
$ws.Range('E20').Value = '  +2.18%  '

$ws.Range('E21').Value = '  +0.88%  '

$ws.Range('D22').Value = "'6.322"
$ws.Range('E22').Value = '  +1.54%  '

$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('D24').Value = "'0.1282"
$ws.Range('E24').Value = '  +5.65%  '

$ws.Range('D25').Value = "'139.45"
$ws.Range('E25').Value = '  -4.73%  '

$ws.Range('D26').Value = "'7.423"
$ws.Range('E26').Value = '  -0.35%  '

$ws.Range('D27').Value = "'16.27"
$ws.Range('E27').Value = '  +2.68%  '

$ws.Range('D28').Value = "'1.432"
$ws.Range('E28').Value = '  +2.29%  '

$ws.Range('D29').Value = "'0.06302"
$ws.Range('E29').Value = '  +6.96%  '

$ws.Range('E30').Value = '  +1.53%  '

$ws.Range('D31').Value = "'3.606"
$ws.Range('E31').Value = '  +5.43%  '

$ws.Range('D32').Value = "'3.441"
$ws.Range('E32').Value = '  +1.16%  '

$ws.Range('E33').Value = '  +2.29%  '

$ws.Range('D34').Value = "'1.013"
$ws.Range('E34').Value = '  +2.30%  '

$ws.Range('D35').Value = "'0.6145"
$ws.Range('E35').Value = '  +8.51%  '

$ws.Range('D36').Value = "'2.424"
$ws.Range('E36').Value = '  +1.26%  '

$ws.Range('E38').Value = '  +0.64%  '

$ws.Range('D39').Value = "'6.099"
$ws.Range('E39').Value = '  +4.34%  '

$ws.Range('D40').Value = "'1.093.19"
$ws.Range('E40').Value = '  +5.84%  '

$ws.Range('D41').Value = "'0.8628"
$ws.Range('E41').Value = '  +0.66%  '

$ws.Range('E42').Value = '  -0.08%  '

$ws.Range('D43').Value = "'100.61"
$ws.Range('E43').Value = '  +0.30%  '

$ws.Range('D44').Value = "'1.823.29"
$ws.Range('E44').Value = '  +1.49%  '

$ws.Range('E45').Value = '  +4.34%  '

$ws.Range('D46').Value = "'58.70"
$ws.Range('E46').Value = '  +4.82%  '

$ws.Range('D47').Value = "'8.200"
$ws.Range('E47').Value = '  +0.98%  '

$ws.Range('D48').Value = "'0.9989"
$ws.Range('E48').Value = '  -0.22%  '

$ws.Range('D49').Value = "'1.512"
$ws.Range('E49').Value = '  +9.19%  '

$ws.Range('D50').Value = "'0.05197"

$ws.Range('D51').Value = "'6.025"
$ws.Range('E51').Value = '  +1.52%  '
